$wb = $excel.ActiveWorkbook

# --- Big_Sheet: update the view selection, it is no longer the active tab ---
$bigSheet = $wb.Worksheets.Item("Big_Sheet")
$bigSheet.Activate()
$bigSheet.Range("B2:B16").Select()

# --- Insert a new "Scores" sheet right before "ExtremeTrees" (2nd tab overall) ---
$extremeTrees = $wb.Worksheets.Item("ExtremeTrees")
$scores = $wb.Worksheets.Add($extremeTrees)
$scores.Name = "Scores"

# Header row
$scores.Range("A1").Value = "model"

# Data rows (model name, score) - write the new string "Multi-Linear Regression"
# before "Score" so the shared-string table ends up in the same order as a
# natural top-to-bottom, left-to-right fill would produce.
$scores.Range("A2").Value = "Extremely Random Trees"
$scores.Range("B2").Value = 0.774
$scores.Range("A3").Value = "Gradient Boost"
$scores.Range("B3").Value = 0.636
$scores.Range("A4").Value = "Random Forest"
$scores.Range("B4").Value = 0.768
$scores.Range("A5").Value = "Multi-Linear Regression"
$scores.Range("B5").Value = 0.175

$scores.Range("B1").Value = "Score"

# Column sizing to fit the longer model names
$scores.Columns.Item(1).ColumnWidth = 22.7109375

# Make Scores the active sheet/tab with its own selection, matching the
# workbook's new "active tab" state
$scores.Activate()
$scores.Range("C12").Select()
